$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$cv = $wb.Worksheets.Item("CONVERTION")

# -----------------------------------------------------------------
# 5. Row 115: clear the PERIOD date, add SL(3-0-0) leave entry,
#    3 days, remark 8/3-5/2023
#    (done first so the new shared string "8/3-5/2023" is created
#    before the other new shared strings below)
# -----------------------------------------------------------------
$ws.Range("A115").ClearContents()
$ws.Range("B115").Value = "SL(3-0-0)"
$ws.Range("H115").Value = 3
$ws.Range("K113").Copy()
$ws.Range("K115").PasteSpecial(-4122)
$ws.Range("K115").Value = "8/3-5/2023"

# -----------------------------------------------------------------
# 2. Row 79: UT(0-0-11) undertime entry
# -----------------------------------------------------------------
$ws.Range("B79").Value = "UT(0-0-11)"
$ws.Range("D79").Value = 0.023000000000000007

# -----------------------------------------------------------------
# 1. Row 78: A(2-0-0) leave entry, 2 days, remark 4/27,28/2022
# -----------------------------------------------------------------
$ws.Range("B78").Value = "A(2-0-0)"
$ws.Range("D78").Value = 2
$ws.Range("K78").Value = "4/27,28/2022"

# -----------------------------------------------------------------
# 3. Row 112: fill in EARNED column (C) with 1.25
# -----------------------------------------------------------------
$ws.Range("C112").Value = 1.25

# -----------------------------------------------------------------
# 4. Row 114: SL(1-0-0) leave entry, 1 day, granted on 45134 (date)
# -----------------------------------------------------------------
$ws.Range("B114").Value = "SL(1-0-0)"
$ws.Range("H114").Value = 1
$ws.Range("K113").Copy()
$ws.Range("K114").PasteSpecial(-4122)
$ws.Range("K114").Value = 45134

# -----------------------------------------------------------------
# 6. Insert a new blank PERIOD row at row 116 (shifts old rows
#    116-144 down to 117-145, keeping the special "last row"
#    border styling on the true last row). The new row takes over
#    the PERIOD date that used to live in row 115.
# -----------------------------------------------------------------
$ws.Rows("116:116").Insert()

$cols = @("A","B","C","D","E","F","G","H","I","J","K")
foreach ($col in $cols) {
    $ws.Range($col + "117").Copy()
    $ws.Range($col + "116").PasteSpecial(-4122)
}
$ws.Range("G116").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),""""," + "Table1[[#This Row],[EARNED]])"
$ws.Range("A116").Value = 45170

# Expand the table to cover the newly inserted row and restore the
# calculated-column formula text on the (new) final row.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K145"))
$ws.Range("G145").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),""""," + "Table1[[#This Row],[EARNED]])"

# -----------------------------------------------------------------
# 7. CONVERTION sheet: update late-calculator inputs (E3/F3), which
#    drives the derived G3 value used by the EARNED formulas above.
# -----------------------------------------------------------------
$cv.Range("E3").Value = 0
$cv.Range("F3").Value = 11

$wb.Save()
